$d = $word.ActiveDocument

# Locate the list paragraph that currently ends with "article added"
# (3rd paragraph overall / 2nd bullet: "This month in WoW article added").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "article added") {
        $target = $p
    }
}

$tr = $target.Range

# Insert the new bullet paragraph right after $target. Collapsing to a point
# that is strictly inside the existing paragraph's text (not touching either
# boundary) makes the engine splice the new paragraph(s) in after the
# current one without disturbing the current paragraph's own content
# (including its bookmark).
$insertionPoint = $d.Range($tr.Start + 1, $tr.Start + 1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Impact" w:hAnsi="Impact" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Major design changes coming to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Starcraft</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> II&#8217;s multiplayer</w:t></w:r></w:p>
<w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($xml)

# Move the _GoBack bookmark from the end of the "article added" paragraph to
# the end of the brand-new paragraph's text, matching the target document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newPara = $d.Paragraphs.Last
$bmRange = $newPara.Range.Duplicate
$bmRange.MoveEnd(1, -1)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
